# Generate Report for Handback
# ------------------------------------------------------------------
# This models a localization "handback" event: the de-de target has
# finished round-tripping and is back in sync, zh-cn's handback file
# is now recorded, and de-de additionally gets a fresh handback
# timestamp. Overview status text + the two status-mirroring columns
# are updated to reflect the new state.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ---- Status text: Overview's per-language status (E2/F2) and each ----
# ---- language sheet's own Status cell (C2) all share one string.  ----
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C2").Value     = "Handed back: in sync with en-US"
$dede.Range("C2").Value     = "Handed back: in sync with en-US"

# ---- zh-cn sheet: Latest Target File (I2) + Latest Handback File (J2) ----
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6835d3a5a5f723b568be90a71d9959554d4d5777/e2e/35fcf230-f3f3-499f-8195-5edfd46dc5d4.md", `
    "", "", "35fcf230-f3f3-499f-8195-5edfd46dc5d4.md")
$zhcn.Range("J2").Value = "35fcf230-f3f3-499f-8195-5edfd46dc5d4.54bad4e69842081a8f0290b0392a202a75f3a2e7.zh-cn.xlf"
# Latest Handback DateTime (K2) for zh-cn was already populated; its
# underlying text now reflects the handback run completion time.
$zhcn.Range("K2").Value = "2016-08-28 22:58:23"

# ---- de-de sheet: Latest Target File (I2) + Latest Handback File (J2) + Latest Handback DateTime (K2) ----
$dede.Hyperlinks.Add($dede.Range("I2"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6835d3a5a5f723b568be90a71d9959554d4d5777/e2e/35fcf230-f3f3-499f-8195-5edfd46dc5d4.md", `
    "", "", "35fcf230-f3f3-499f-8195-5edfd46dc5d4.md")
$dede.Range("J2").Value = "35fcf230-f3f3-499f-8195-5edfd46dc5d4.54bad4e69842081a8f0290b0392a202a75f3a2e7.de-de.xlf"
$dede.Range("K2").Value = "2016-08-28 22:58:30"

# ---- Column widths: widen the columns that now hold longer content ----
# (The workbook author re-ran AutoFit after filling these columns in;
# reproduce the resulting widths as closely as this host allows.)
$overview.Columns.Item(5).ColumnWidth = 29.144371396019366
$overview.Columns.Item(6).ColumnWidth = 29.144371396019366

$zhcn.Columns.Item(3).ColumnWidth  = 29.144371396019366
$zhcn.Columns.Item(9).ColumnWidth  = 38.79187157040546
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

$dede.Columns.Item(3).ColumnWidth  = 29.144371396019366
$dede.Columns.Item(9).ColumnWidth  = 38.79187157040546
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664

